# This script reproduces the authored change:
#  - A4 changes from the literal text "2025-11-30" to a live =TODAY() formula,
#    formatted with the built-in short-date number format (numFmtId 14). This
#    also drops the now-unused "2025-11-30" shared string and re-numbers the
#    shared-string indices used further down the sheet (handled automatically
#    by the workbook when the string is no longer referenced).
#  - The active selection on the "Tests" sheet moves from D9 to D5.
#  - The workbook window position/size is updated to the author's last-saved
#    window geometry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A4: replace the static date text with a volatile =TODAY() formula ---
# Apply the number format first so Excel records the built-in date format
# (numFmtId 14) instead of auto-generating a custom one when the formula is
# entered.
$ws.Range("A4").NumberFormat = "mm-dd-yy"
$ws.Range("A4").Formula = "=TODAY()"

# --- Update the active cell/selection shown when the sheet is reopened ---
$ws.Range("D5").Select() | Out-Null

# --- Restore the author's last window position/size ---
$win = $excel.ActiveWindow
$win.Left = 23025
$win.Top = 3600
$win.Width = 23145
$win.Height = 14685
